# Update the Sweden Summary country-indicator figures with more precise
# (two-decimal) values, per the "Update country data files" commit.
#
# The target cells hold their numeric-looking values as literal TEXT
# (shared strings), matching how the source workbook stores every
# indicator. Assigning a numeric-looking string straight to .Value /
# .Value2 / .Formula makes Excel auto-convert it to a real number, so we
# briefly force Text number-formatting while we type the new values, then
# restore each cell's original formatting by copying it back from an
# untouched cell that already carries the same (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (old value shown in the comment for reference).
$updates = @{
    "B11" = "107.92"  # was 107.9
    "C11" = "1.67"    # was 1.7
    "D11" = "109.59"  # was 109.6
    "B12" = "33.16"   # was 33.2
    "C12" = "30.94"   # was 30.9
    "D12" = "64.09"   # was 64.1
    "C33" = "3.63"    # was 3.6
    "D33" = "68.83"   # was 68.8
    "B34" = "25.93"   # was 25.9
    "C34" = "39.72"   # was 39.7
    "D34" = "65.65"   # was 65.6
    "B36" = "94.58"   # was 94.6
    "C36" = "5.26"    # was 5.3
    "D36" = "99.85"   # was 99.8
    "B40" = "21.69"   # was 21.7
    "C40" = "37.81"   # was 37.8
    "D40" = "59.51"   # was 59.5
}

# A cell that is untouched by this edit and already uses the plain/default
# style shared by every cell we are about to rewrite (used below to restore
# formatting after the temporary Text number-format trick).
$styleDonor = $ws.Range("A22")

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]

    $styleDonor.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false
